$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "D2" "42.367.92"
Set-TextValue "D3" "2.219.57"
Set-TextValue "E3" "  -2.07%  "
Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.40%  "
Set-TextValue "D5" "107.89"
Set-TextValue "E5" "  -9.64%  "
Set-TextValue "D6" "296.80"
Set-TextValue "E6" "  +11.96%  "
Set-TextValue "E7" "  -3.15%  "
Set-TextValue "E8" "  -0.12%  "
Set-TextValue "E9" "  -3.08%  "
Set-TextValue "D10" "43.61"
Set-TextValue "E10" "  -8.09%  "
Set-TextValue "E11" "  -3.21%  "
Set-TextValue "D12" "54.49"
Set-TextValue "E12" "  +0.59%  "
Set-TextValue "D13" "8.74"
Set-TextValue "E13" "  -4.14%  "
Set-TextValue "D14" "1.01"
Set-TextValue "E14" "  +11.03%  "
Set-TextValue "E15" "  -2.99%  "
Set-TextValue "D16" "15.03"
Set-TextValue "E16" "  -1.95%  "
Set-TextValue "D17" "2.550.94"
Set-TextValue "E17" "  -2.08%  "
Set-TextValue "D18" "2.228.82"
Set-TextValue "E18" "  -1.80%  "
Set-TextValue "D19" "42.284.36"
Set-TextValue "E19" "  -2.79%  "
Set-TextValue "D20" "7.38"
Set-TextValue "E20" "  +7.50%  "
Set-TextValue "E21" "  -4.12%  "
Set-TextValue "D22" "72.38"
Set-TextValue "E22" "  +0.46%  "
Set-TextValue "D23" "3.47"
Set-TextValue "E23" "  +21.36%  "
Set-TextValue "D24" "2.30"
Set-TextValue "E24" "  -3.58%  "
Set-TextValue "D25" "228.46"
Set-TextValue "E25" "  -2.91%  "
Set-TextValue "D26" "9.10"
Set-TextValue "E26" "  -4.28%  "
Set-TextValue "E27" "  -1.72%  "
Set-TextValue "E28" "  -2.55%  "
Set-TextValue "E29" "  -0.89%  "
Set-TextValue "D30" "38.10"
Set-TextValue "E30" "  -8.17%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D31" "174.06"
Set-TextValue "E31" "  +1.30%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D32" "3.19"
Set-TextValue "E32" "  -5.12%  "
Set-TextValue "D33" "20.95"
Set-TextValue "E33" "  -3.54%  "
Set-TextValue "D34" "0.0897"
Set-TextValue "E34" "  -1.67%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D35" "5.10"
Set-TextValue "E35" "  +12.19%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D36" "5.59"
Set-TextValue "E36" "  -1.74%  "
Set-TextValue "D37" "4.34"
Set-TextValue "E37" "  +0.63%  "
Set-TextValue "E38" "  -3.14%  "
Set-TextValue "D39" "0.0376"
Set-TextValue "E39" "  -0.03%  "
Set-TextValue "E40" "  -3.84%  "
Set-TextValue "D41" "2.41"
Set-TextValue "E41" "  -5.33%  "
Set-TextValue "D42" "72.11"
Set-TextValue "E42" "  -2.87%  "
Set-TextValue "E43" "  -2.04%  "
Set-TextValue "E44" "  +0.16%  "
Set-TextValue "D45" "12.61"
Set-TextValue "E45" "  -9.45%  "
Set-TextValue "E46" "  -4.70%  "
Set-TextValue "D47" "5.40"
Set-TextValue "E47" "  -6.19%  "
Set-TextValue "E48" "  +5.03%  "
Set-TextValue "D49" "103.26"
Set-TextValue "E49" "  +1.87%  "
Set-TextValue "D50" "1.66"
Set-TextValue "E50" "  +6.73%  "
Set-TextValue "E51" "  -1.30%  "
